# Timesheet update for 18th March.
# - Insert a new task row ("Java fx: Editscreen") above the "Mail class:" row,
#   logging 6 hours on 18 Mar (column X).
# - Log an extra 6 hours on 18 Mar for four other in-progress tasks that were
#   already logged on 17 Mar (column W): "DB: Implementation of database in
#   java.", "KTN: Functionality for GUI class", "Java fx: Main screen" and
#   "Java fx: Viewscreen and logic".
# - Move the active selection to reflect the new bottom of the log (W46).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 33 ("Mail class:"), pushing everything below
# (Mail class:, the 357-total row, and the notes rows) down by one. Excel
# auto-extends the shared SUM formulas in row 3 (C6:C39 -> C6:C40, etc.) and
# the sheet dimension as part of this insert.
$ws.Rows.Item(33).Insert()

# New task entry in the freshly inserted row 33.
$ws.Range("A33").Value = "Java fx: Editscreen"
$ws.Range("X33").Value = 6

# Extra hours logged for 18 Mar (column X) on tasks already worked on 17 Mar.
$ws.Range("X26").Value = 6
$ws.Range("X30").Value = 6
$ws.Range("X31").Value = 6
$ws.Range("X32").Value = 6

# Reflect the new editing position.
$ws.Range("W46").Select()
